$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, shifting existing rows 100-117 down to 101-118.
$ws.Rows.Item(100).EntireRow.Insert()

# Populate the newly inserted row 100 with the new weekly price record.
$ws.Cells.Item(100, 1).Value = 7
$ws.Cells.Item(100, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(100, 3).Value = "Ñuble"
$ws.Cells.Item(100, 4).Value = 44995
$ws.Cells.Item(100, 5).Value = 16
$ws.Cells.Item(100, 6).Value = 100112037
$ws.Cells.Item(100, 7).Value = "Cebollín"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 50
$ws.Cells.Item(100, 11).Value = 6000
$ws.Cells.Item(100, 12).Value = 6000
$ws.Cells.Item(100, 13).Value = 6000
$ws.Cells.Item(100, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(100, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(100, 16).Value = 167
$ws.Cells.Item(100, 17).Value = 36
$ws.Cells.Item(100, 18).Value = "Hortaliza"
